$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3916.75
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 4000.087
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4000.087
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -4350.087

$ws.Range("H44").Value = 36623.332
$ws.Range("J44").Value = 36623.332
$ws.Range("L44").Value = 36623.332
$ws.Range("N44").Value = -37547.332

$ws.Range("H129").Value = 1671.0714
$ws.Range("I129").Value = 1232.5
$ws.Range("K129").Value = 3697.5
$ws.Range("M129").Value = 1302.5

$ws.Range("H137").Value = 2585.7917
$ws.Range("I137").Value = 1968.8
$ws.Range("K137").Value = 5906.4
$ws.Range("M137").Value = -3356.4

$ws.Range("H138").Value = 1768.0103
$ws.Range("J138").Value = 2478.8772
$ws.Range("L138").Value = 7436.6316
$ws.Range("N138").Value = -17716.6316


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2233.7058
$ws.Range("I2").Value = 2134.0833
$ws.Range("K2").Value = 2134.0833
$ws.Range("M2").Value = -2021.0833

$ws.Range("H32").Value = 29259304
$ws.Range("I32").Value = 41111300
$ws.Range("K32").Value = 41111300
$ws.Range("M32").Value = -41111013

$ws.Range("H39").Value = 5278.75
$ws.Range("I39").Value = 4705
$ws.Range("J39").Value = 7000
$ws.Range("K39").Value = 4705
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = -4185
$ws.Range("N39").Value = -8040

$ws.Range("H61").Value = 4238.222
$ws.Range("I61").Value = 4614.5
$ws.Range("K61").Value = 4614.5
$ws.Range("M61").Value = -4402.5

$ws.Range("H63").Value = 3996.6365
$ws.Range("I63").Value = 2366.125
$ws.Range("J63").Value = 4928.357
$ws.Range("K63").Value = 2366.125
$ws.Range("L63").Value = 4928.357
$ws.Range("M63").Value = -1680.125
$ws.Range("N63").Value = -6300.357

$ws.Range("H66").Value = 3996.6365
$ws.Range("I66").Value = 2366.125
$ws.Range("J66").Value = 4928.357
$ws.Range("K66").Value = 11830.625
$ws.Range("L66").Value = 24641.785
$ws.Range("M66").Value = -8398.625
$ws.Range("N66").Value = -31505.785

$ws.Range("H116").Value = 2233.7058
$ws.Range("I116").Value = 2134.0833
$ws.Range("K116").Value = 2134.0833
$ws.Range("M116").Value = 159.9167000000002

$ws.Range("H136").Value = 4238.222
$ws.Range("I136").Value = 4614.5
$ws.Range("K136").Value = 13843.5
$ws.Range("M136").Value = -11293.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2233.7058
$ws.Range("I3").Value = 2134.0833
$ws.Range("K3").Value = 2134.0833
$ws.Range("M3").Value = -2020.0833

$ws.Range("H56").Value = 22100
$ws.Range("I56").Value = 4100
$ws.Range("J56").Value = 23900
$ws.Range("K56").Value = 4100
$ws.Range("L56").Value = 23900
$ws.Range("M56").Value = -3361
$ws.Range("N56").Value = -25378

$ws.Range("H134").Value = 3763895.2
$ws.Range("I134").Value = 5954768
$ws.Range("K134").Value = 17864304
$ws.Range("M134").Value = -17861769


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4013.0386
$ws.Range("J31").Value = 7692.8276
$ws.Range("L31").Value = 7692.8276
$ws.Range("N31").Value = -8282.827600000001

$ws.Range("H34").Value = 4013.0386
$ws.Range("J34").Value = 7692.8276
$ws.Range("L34").Value = 7692.8276
$ws.Range("N34").Value = -8096.8276

$ws.Range("H58").Value = 2569.1
$ws.Range("I58").Value = 2450.7144
$ws.Range("K58").Value = 2450.7144
$ws.Range("M58").Value = -2247.7144

$ws.Range("H136").Value = 2569.1
$ws.Range("I136").Value = 2450.7144
$ws.Range("K136").Value = 7352.1432
$ws.Range("M136").Value = -4802.1432


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2441.8572
$ws.Range("I33").Value = 2219.25
$ws.Range("J33").Value = 2738.6667
$ws.Range("K33").Value = 13315.5
$ws.Range("L33").Value = 16432.0002
$ws.Range("M33").Value = -13032.5
$ws.Range("N33").Value = -16998.0002


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 2525.4119
$ws.Range("I102").Value = 2328.8
$ws.Range("K102").Value = 2328.8
$ws.Range("M102").Value = -706.8000000000002

$ws.Range("H122").Value = 1755
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1755
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5265
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10165


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3702.6365
$ws.Range("I7").Value = 3616.25
$ws.Range("K7").Value = 3616.25
$ws.Range("M7").Value = -3504.25

$ws.Range("H126").Value = 3702.6365
$ws.Range("I126").Value = 3616.25
$ws.Range("K126").Value = 10848.75
$ws.Range("M126").Value = -8378.75

$ws.Range("H132").Value = 4624.5
$ws.Range("I132").Value = 4199.4
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 12598.2
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -10068.2
$ws.Range("N132").Value = -21059


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5505.7144
$ws.Range("I62").Value = 4280
$ws.Range("J62").Value = 6425
$ws.Range("K62").Value = 4280
$ws.Range("L62").Value = 6425
$ws.Range("M62").Value = -3656
$ws.Range("N62").Value = -7673

$ws.Range("H65").Value = 5505.7144
$ws.Range("I65").Value = 4280
$ws.Range("J65").Value = 6425
$ws.Range("K65").Value = 21400
$ws.Range("L65").Value = 32125
$ws.Range("M65").Value = -18280
$ws.Range("N65").Value = -38365

$ws.Range("H126").Value = 5918.3
$ws.Range("I126").Value = 7862.1665
$ws.Range("J126").Value = 3002.5
$ws.Range("K126").Value = 23586.4995
$ws.Range("L126").Value = 9007.5
$ws.Range("M126").Value = -21116.4995
$ws.Range("N126").Value = -13947.5

$ws.Range("H132").Value = 2483.7222
$ws.Range("I132").Value = 2306.2942
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 6918.882599999999
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -4388.882599999999
$ws.Range("N132").Value = -21560
